$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (k values)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary block, rows 14-17, columns A (label) and B (value)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the desired formatting (bold, size 12, vertically centered) on a
# scratch cell first, so it collapses into a single new style entry, then
# copy that format onto the summary block in one shot.
$scratch = $ws.Range("AA1")
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108

$target = $ws.Range("A14:B17")
[void]$scratch.Copy()
[void]$target.PasteSpecial(-4122)  # xlPasteFormats
[void]$scratch.Clear()

# Update the selection to match the edited range
[void]$ws.Range("A14:B17").Select()

# Page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

[void]$wb.Save()
